$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Camote, "1a (cosecha)", date 2022-05-30) is
# inserted above the existing row 209, pushing rows 209:253 down to 210:254.
$ws.Rows.Item(209).Insert()

$ws.Range("A209").Value = 5
$ws.Range("B209").Value = "Macroferia Regional de Talca"
$ws.Range("C209").Value = "Maule"
$ws.Range("D209").Value = 44711
$ws.Range("E209").Value = 7
$ws.Range("F209").Value = 100112045
$ws.Range("G209").Value = "Zapallo"
$ws.Range("H209").Value = "Camote"
$ws.Range("I209").Value = "1a (cosecha)"
$ws.Range("J209").Value = 900
$ws.Range("K209").Value = 350
$ws.Range("L209").Value = 350
$ws.Range("M209").Value = 350
$ws.Range("N209").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O209").Value = "Región del Maule"
$ws.Range("P209").Value = 350
$ws.Range("Q209").Value = 1
$ws.Range("R209").Value = "Hortaliza"
